$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename "Sheet1" -> "Data" (the _FilterDatabase defined name that
#    points at Sheet1!$A$1:$G$1001 is updated automatically by the
#    rename).
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Sheet1")
$wsData.Name = "Data"

# ---------------------------------------------------------------------
# 2) Rebuild the "Legend" sheet.
#    It used to hold volatile dynamic-array SORT(UNIQUE(...)) helper
#    formulas spilling into A1:F6 plus a separate list of legend labels
#    in A9:A14. Replace all of that with a small, static two-column
#    legend (score -> label) and turn it into a proper table.
# ---------------------------------------------------------------------
$wsLegend = $wb.Worksheets.Item("Legend")

# Break the spilled dynamic-array formulas first (clearing the spilled
# cells directly raises "You cannot change part of an array"; clearing
# the formula on the anchor cells removes the whole spill cleanly).
$wsLegend.Range("A1").Formula = ""
$wsLegend.Range("B1").Formula = ""
$wsLegend.Range("C1").Formula = ""
$wsLegend.Range("D1").Formula = ""
$wsLegend.Range("E1").Formula = ""
$wsLegend.Range("F1").Formula = ""

# Drop whatever is left over (the old 1-5 helper numbers and the old
# A9:A14 legend text list).
$wsLegend.Range("A2:F14").ClearContents()

# Write the new legend body first (rows 2-7) so the shared-string table
# picks up the label text right after the existing strings, then write
# the header row last (it reuses two brand-new "Column1"/"Column2"
# strings that must land at the very end of the shared-string table).
$wsLegend.Range("A2").Value = 0
$wsLegend.Range("B2").Value = " No Concern "
$wsLegend.Range("A3").Value = 1
$wsLegend.Range("B3").Value = " Minimal Concern "
$wsLegend.Range("A4").Value = 2
$wsLegend.Range("B4").Value = " Slightly Concerning "
$wsLegend.Range("A5").Value = 3
$wsLegend.Range("B5").Value = " Moderately Concerning "
$wsLegend.Range("A6").Value = 4
$wsLegend.Range("B6").Value = " Highly Concerning "
$wsLegend.Range("A7").Value = 5
$wsLegend.Range("B7").Value = " Extremely Concerning "

$wsLegend.Range("A1").Value = "Column1"
$wsLegend.Range("B1").Value = "Column2"

# Turn A1:B7 into a real table ("Table1") with an autofilter, matching
# the new xl/tables/table1.xml part.
$tbl = $wsLegend.ListObjects.Add(1, $wsLegend.Range("A1:B7"), 0, 1)
$tbl.Name = "Table1"

# Give the two columns a bit of explicit width, and select the table
# range on the Legend sheet (cosmetic, matches the saved sheet view).
$wsLegend.Columns("A:B").ColumnWidth = 10.43
$wsLegend.Activate()
[void]$wsLegend.Range("A1:B7").Select()

# Restore "Data" as the active/selected sheet, same as before the edit.
[void]$wsData.Activate()
